$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy formatting from the last existing data row (200) down into the new rows
# (201-210) so the new cells reuse the same cell style ("import") as every other
# data row, instead of picking up the plain column default style.
$ws.Range("A200:C200").Copy()
$ws.Range("A201:C210").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Introduce the new shared strings in the same order the original authoring
# session created them in (so the shared-string table lines up cell-for-cell).
$ws.Range("C205").Value = 'Seznam buildů'
$ws.Range("C202").Value = 'V této sekci je možné spravovat a vytvářet buildy.'
$ws.Range("B203").Value = 'lab.build.index.menu'
$ws.Range("B204").Value = 'lab.build.create.menu'
$ws.Range("C204").Value = 'Nový build'
$ws.Range("B205").Value = 'lab.build.list.menu'
$ws.Range("B201").Value = 'lab.build.title'
$ws.Range("B202").Value = 'lab.build.subtitle'
$ws.Range("B206").Value = 'lab.build.create.title'
$ws.Range("B207").Value = 'lab.build.create.subtitle'
$ws.Range("B208").Value = 'lab.build.name.label'
$ws.Range("C208").Value = 'Jméno'
$ws.Range("B209").Value = 'lab.build.name.label.tooltip'
$ws.Range("C209").Value = 'Jméno buildu musí být unikátní; později bude možné ho použít pro evidenci požitků (vapování).'
$ws.Range("B210").Value = 'lab.build.atomizerId.label'

# Remaining cells only reuse strings that already exist in the workbook, so
# their write order does not affect the shared-string table.
$ws.Range("A201").Value = 'cs'
$ws.Range("A202").Value = 'cs'
$ws.Range("A203").Value = 'cs'
$ws.Range("A204").Value = 'cs'
$ws.Range("A205").Value = 'cs'
$ws.Range("A206").Value = 'cs'
$ws.Range("A207").Value = 'cs'
$ws.Range("A208").Value = 'cs'
$ws.Range("A209").Value = 'cs'
$ws.Range("A210").Value = 'cs'
$ws.Range("C201").Value = 'Buildy'
$ws.Range("C203").Value = 'Obecné'
$ws.Range("C206").Value = 'Nový build'
$ws.Range("C207").Value = ' '
$ws.Range("C210").Value = 'Atomizér'

# Update the visible selection to match the post-edit cursor position.
$ws.Range("B196").Select()
